$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dimension by touching the new extent; Q/R cleanup for existing rows ---

# Rows where Q and R both reset to 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("R69").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = 0

# Rows where only R resets to 0
$ws.Range("R91").Value = 0
$ws.Range("R99").Value = 0
$ws.Range("R104").Value = 0
$ws.Range("R111").Value = 0
$ws.Range("R119").Value = 0
$ws.Range("R144").Value = 0
$ws.Range("R150").Value = 0
$ws.Range("R165").Value = 0
$ws.Range("R167").Value = 0
$ws.Range("R176").Value = 0
$ws.Range("R181").Value = 0
$ws.Range("R187").Value = 0
$ws.Range("R214").Value = 0
$ws.Range("R247").Value = 0
$ws.Range("R261").Value = 0
$ws.Range("R266").Value = 0
$ws.Range("R277").Value = 0
$ws.Range("R283").Value = 0
$ws.Range("R292").Value = 0
$ws.Range("R307").Value = 0
$ws.Range("R320").Value = 0
$ws.Range("R334").Value = 0
$ws.Range("R342").Value = 0
$ws.Range("R350").Value = 0
$ws.Range("R358").Value = 0
$ws.Range("R368").Value = 0

# Row 377: O updated
$ws.Range("O377").Value = 2

# Row 380: O updated, R filled in as 0 (was blank)
$ws.Range("O380").Value = 1
$ws.Range("R380").Value = 0

# Row 381: R filled in as 0 (was blank)
$ws.Range("R381").Value = 0

# --- Append new weekly rows 382-411 ---
$dateFmt = $ws.Range("A2").NumberFormat

$ws.Range("A382").Value = 45474
$ws.Range("A382").NumberFormat = $dateFmt
$ws.Range("B382").Value = 4725.10009765625
$ws.Range("C382").Value = 4871
$ws.Range("D382").Value = 4612.5
$ws.Range("E382").Value = 4853.10009765625
$ws.Range("G382").Value = 2478739
$ws.Range("H382").Value = 2024
$ws.Range("I382").Value = 7
$ws.Range("J382").Value = 1
$ws.Range("K382").Value = 0
$ws.Range("L382").Value = 0
$ws.Range("M382").Value = 0
$ws.Range("N382").Value = 27
$ws.Range("O382").Value = 0
$ws.Range("P382").Value = 0
$ws.Range("Q382").Value = 0

$ws.Range("A383").Value = 45481
$ws.Range("A383").NumberFormat = $dateFmt
$ws.Range("B383").Value = 4860
$ws.Range("C383").Value = 5016.9501953125
$ws.Range("D383").Value = 4731.10009765625
$ws.Range("E383").Value = 4943.64990234375
$ws.Range("G383").Value = 2578794
$ws.Range("H383").Value = 2024
$ws.Range("I383").Value = 7
$ws.Range("J383").Value = 8
$ws.Range("K383").Value = 0
$ws.Range("L383").Value = 0
$ws.Range("M383").Value = 0
$ws.Range("N383").Value = 28
$ws.Range("O383").Value = 0
$ws.Range("P383").Value = 0
$ws.Range("Q383").Value = 1

$ws.Range("A384").Value = 45488
$ws.Range("A384").NumberFormat = $dateFmt
$ws.Range("B384").Value = 5126.4501953125
$ws.Range("C384").Value = 5164
$ws.Range("D384").Value = 4943.64990234375
$ws.Range("E384").Value = 5010.7001953125
$ws.Range("G384").Value = 3750963
$ws.Range("H384").Value = 2024
$ws.Range("I384").Value = 7
$ws.Range("J384").Value = 15
$ws.Range("K384").Value = 0
$ws.Range("L384").Value = 0
$ws.Range("M384").Value = 0
$ws.Range("N384").Value = 29
$ws.Range("O384").Value = 0
$ws.Range("P384").Value = 0
$ws.Range("Q384").Value = 0

$ws.Range("A385").Value = 45495
$ws.Range("A385").NumberFormat = $dateFmt
$ws.Range("B385").Value = 4955
$ws.Range("C385").Value = 5205
$ws.Range("D385").Value = 4922.75
$ws.Range("E385").Value = 5071.60009765625
$ws.Range("G385").Value = 2249183
$ws.Range("H385").Value = 2024
$ws.Range("I385").Value = 7
$ws.Range("J385").Value = 22
$ws.Range("K385").Value = 0
$ws.Range("L385").Value = 0
$ws.Range("M385").Value = 0
$ws.Range("N385").Value = 30
$ws.Range("O385").Value = 0
$ws.Range("P385").Value = 0
$ws.Range("Q385").Value = 0

$ws.Range("A386").Value = 45502
$ws.Range("A386").NumberFormat = $dateFmt
$ws.Range("B386").Value = 5144
$ws.Range("C386").Value = 5175.5498046875
$ws.Range("D386").Value = 4883
$ws.Range("E386").Value = 4909.7001953125
$ws.Range("G386").Value = 2107424
$ws.Range("H386").Value = 2024
$ws.Range("I386").Value = 7
$ws.Range("J386").Value = 29
$ws.Range("K386").Value = 0
$ws.Range("L386").Value = 0
$ws.Range("M386").Value = 0
$ws.Range("N386").Value = 31
$ws.Range("O386").Value = 0
$ws.Range("P386").Value = 0
$ws.Range("Q386").Value = 2

$ws.Range("A387").Value = 45509
$ws.Range("A387").NumberFormat = $dateFmt
$ws.Range("B387").Value = 4861.5
$ws.Range("C387").Value = 5063.89990234375
$ws.Range("D387").Value = 4801
$ws.Range("E387").Value = 4989.9501953125
$ws.Range("G387").Value = 2286673
$ws.Range("H387").Value = 2024
$ws.Range("I387").Value = 8
$ws.Range("J387").Value = 5
$ws.Range("K387").Value = 0
$ws.Range("L387").Value = 0
$ws.Range("M387").Value = 0
$ws.Range("N387").Value = 32
$ws.Range("O387").Value = 0
$ws.Range("P387").Value = 0
$ws.Range("Q387").Value = 0

$ws.Range("A388").Value = 45516
$ws.Range("A388").NumberFormat = $dateFmt
$ws.Range("B388").Value = 4999
$ws.Range("C388").Value = 5048
$ws.Range("D388").Value = 4902.7001953125
$ws.Range("E388").Value = 5021.2998046875
$ws.Range("G388").Value = 782256
$ws.Range("H388").Value = 2024
$ws.Range("I388").Value = 8
$ws.Range("J388").Value = 12
$ws.Range("K388").Value = 0
$ws.Range("L388").Value = 0
$ws.Range("M388").Value = 0
$ws.Range("N388").Value = 33
$ws.Range("O388").Value = 0
$ws.Range("P388").Value = 0
$ws.Range("Q388").Value = 0

$ws.Range("A389").Value = 45523
$ws.Range("A389").NumberFormat = $dateFmt
$ws.Range("B389").Value = 5040
$ws.Range("C389").Value = 5164
$ws.Range("D389").Value = 4880
$ws.Range("E389").Value = 4901.5
$ws.Range("G389").Value = 2444554
$ws.Range("H389").Value = 2024
$ws.Range("I389").Value = 8
$ws.Range("J389").Value = 19
$ws.Range("K389").Value = 0
$ws.Range("L389").Value = 0
$ws.Range("M389").Value = 0
$ws.Range("N389").Value = 34
$ws.Range("O389").Value = 0
$ws.Range("P389").Value = 0
$ws.Range("Q389").Value = 0

$ws.Range("A390").Value = 45530
$ws.Range("A390").NumberFormat = $dateFmt
$ws.Range("B390").Value = 4903
$ws.Range("C390").Value = 5155.5498046875
$ws.Range("D390").Value = 4900
$ws.Range("E390").Value = 4927.4501953125
$ws.Range("G390").Value = 2628210
$ws.Range("H390").Value = 2024
$ws.Range("I390").Value = 8
$ws.Range("J390").Value = 26
$ws.Range("K390").Value = 0
$ws.Range("L390").Value = 0
$ws.Range("M390").Value = 0
$ws.Range("N390").Value = 35
$ws.Range("O390").Value = 0
$ws.Range("P390").Value = 0
$ws.Range("Q390").Value = 0

$ws.Range("A391").Value = 45537
$ws.Range("A391").NumberFormat = $dateFmt
$ws.Range("B391").Value = 4977
$ws.Range("C391").Value = 5361.14990234375
$ws.Range("D391").Value = 4951
$ws.Range("E391").Value = 5303.4501953125
$ws.Range("G391").Value = 2877063
$ws.Range("H391").Value = 2024
$ws.Range("I391").Value = 9
$ws.Range("J391").Value = 2
$ws.Range("K391").Value = 0
$ws.Range("L391").Value = 0
$ws.Range("M391").Value = 0
$ws.Range("N391").Value = 36
$ws.Range("O391").Value = 0
$ws.Range("P391").Value = 0
$ws.Range("Q391").Value = 0

$ws.Range("A392").Value = 45544
$ws.Range("A392").NumberFormat = $dateFmt
$ws.Range("B392").Value = 5303.4501953125
$ws.Range("C392").Value = 5449
$ws.Range("D392").Value = 5155
$ws.Range("E392").Value = 5187.0498046875
$ws.Range("G392").Value = 2233777
$ws.Range("H392").Value = 2024
$ws.Range("I392").Value = 9
$ws.Range("J392").Value = 9
$ws.Range("K392").Value = 0
$ws.Range("L392").Value = 0
$ws.Range("M392").Value = 0
$ws.Range("N392").Value = 37
$ws.Range("O392").Value = 0
$ws.Range("P392").Value = 0
$ws.Range("Q392").Value = 0

$ws.Range("A393").Value = 45551
$ws.Range("A393").NumberFormat = $dateFmt
$ws.Range("B393").Value = 5188.0498046875
$ws.Range("C393").Value = 5443.9501953125
$ws.Range("D393").Value = 5176.14990234375
$ws.Range("E393").Value = 5320.5498046875
$ws.Range("G393").Value = 1547072
$ws.Range("H393").Value = 2024
$ws.Range("I393").Value = 9
$ws.Range("J393").Value = 16
$ws.Range("K393").Value = 0
$ws.Range("L393").Value = 0
$ws.Range("M393").Value = 0
$ws.Range("N393").Value = 38
$ws.Range("O393").Value = 0
$ws.Range("P393").Value = 0
$ws.Range("Q393").Value = 0

$ws.Range("A394").Value = 45558
$ws.Range("A394").NumberFormat = $dateFmt
$ws.Range("B394").Value = 5355.14990234375
$ws.Range("C394").Value = 5484.85009765625
$ws.Range("D394").Value = 5090
$ws.Range("E394").Value = 5102.2998046875
$ws.Range("G394").Value = 2320945
$ws.Range("H394").Value = 2024
$ws.Range("I394").Value = 9
$ws.Range("J394").Value = 23
$ws.Range("K394").Value = 0
$ws.Range("L394").Value = 0
$ws.Range("M394").Value = 0
$ws.Range("N394").Value = 39
$ws.Range("O394").Value = 1
$ws.Range("P394").Value = 0
$ws.Range("Q394").Value = 0

$ws.Range("A395").Value = 45565
$ws.Range("A395").NumberFormat = $dateFmt
$ws.Range("B395").Value = 5102.2998046875
$ws.Range("C395").Value = 5215
$ws.Range("D395").Value = 4696.10009765625
$ws.Range("E395").Value = 4737.5498046875
$ws.Range("G395").Value = 4598805
$ws.Range("H395").Value = 2024
$ws.Range("I395").Value = 9
$ws.Range("J395").Value = 30
$ws.Range("K395").Value = 0
$ws.Range("L395").Value = 0
$ws.Range("M395").Value = 0
$ws.Range("N395").Value = 40
$ws.Range("O395").Value = 0
$ws.Range("P395").Value = 0
$ws.Range("Q395").Value = 0

$ws.Range("A396").Value = 45572
$ws.Range("A396").NumberFormat = $dateFmt
$ws.Range("B396").Value = 4741
$ws.Range("C396").Value = 4752.89990234375
$ws.Range("D396").Value = 4435
$ws.Range("E396").Value = 4572.7001953125
$ws.Range("G396").Value = 4661035
$ws.Range("H396").Value = 2024
$ws.Range("I396").Value = 10
$ws.Range("J396").Value = 7
$ws.Range("K396").Value = 0
$ws.Range("L396").Value = 0
$ws.Range("M396").Value = 0
$ws.Range("N396").Value = 41
$ws.Range("O396").Value = 0
$ws.Range("P396").Value = 0
$ws.Range("Q396").Value = 0

$ws.Range("A397").Value = 45579
$ws.Range("A397").NumberFormat = $dateFmt
$ws.Range("B397").Value = 4204
$ws.Range("C397").Value = 4299
$ws.Range("D397").Value = 3975.10009765625
$ws.Range("E397").Value = 3986.699951171875
$ws.Range("G397").Value = 8747513
$ws.Range("H397").Value = 2024
$ws.Range("I397").Value = 10
$ws.Range("J397").Value = 14
$ws.Range("K397").Value = 0
$ws.Range("L397").Value = 0
$ws.Range("M397").Value = 0
$ws.Range("N397").Value = 42
$ws.Range("O397").Value = 0
$ws.Range("P397").Value = 0
$ws.Range("Q397").Value = 0

$ws.Range("A398").Value = 45586
$ws.Range("A398").NumberFormat = $dateFmt
$ws.Range("B398").Value = 4000.050048828125
$ws.Range("C398").Value = 4224.7998046875
$ws.Range("D398").Value = 3961
$ws.Range("E398").Value = 4052.199951171875
$ws.Range("G398").Value = 2738848
$ws.Range("H398").Value = 2024
$ws.Range("I398").Value = 10
$ws.Range("J398").Value = 21
$ws.Range("K398").Value = 0
$ws.Range("L398").Value = 0
$ws.Range("M398").Value = 0
$ws.Range("N398").Value = 43
$ws.Range("O398").Value = 0
$ws.Range("P398").Value = 0
$ws.Range("Q398").Value = 0

$ws.Range("A399").Value = 45593
$ws.Range("A399").NumberFormat = $dateFmt
$ws.Range("B399").Value = 4052.199951171875
$ws.Range("C399").Value = 4094.300048828125
$ws.Range("D399").Value = 3876
$ws.Range("E399").Value = 4001.60009765625
$ws.Range("G399").Value = 1587934
$ws.Range("H399").Value = 2024
$ws.Range("I399").Value = 10
$ws.Range("J399").Value = 28
$ws.Range("K399").Value = 0
$ws.Range("L399").Value = 0
$ws.Range("M399").Value = 0
$ws.Range("N399").Value = 44
$ws.Range("O399").Value = 0
$ws.Range("P399").Value = 0
$ws.Range("Q399").Value = 0

$ws.Range("A400").Value = 45600
$ws.Range("A400").NumberFormat = $dateFmt
$ws.Range("B400").Value = 4000.75
$ws.Range("C400").Value = 4006.10009765625
$ws.Range("D400").Value = 3860.14990234375
$ws.Range("E400").Value = 3874.5
$ws.Range("G400").Value = 1810975
$ws.Range("H400").Value = 2024
$ws.Range("I400").Value = 11
$ws.Range("J400").Value = 4
$ws.Range("K400").Value = 0
$ws.Range("L400").Value = 0
$ws.Range("M400").Value = 0
$ws.Range("N400").Value = 45
$ws.Range("O400").Value = 0
$ws.Range("P400").Value = 0
$ws.Range("Q400").Value = 2

$ws.Range("A401").Value = 45607
$ws.Range("A401").NumberFormat = $dateFmt
$ws.Range("B401").Value = 3899.64990234375
$ws.Range("C401").Value = 3940
$ws.Range("D401").Value = 3751
$ws.Range("E401").Value = 3823.85009765625
$ws.Range("G401").Value = 1827153
$ws.Range("H401").Value = 2024
$ws.Range("I401").Value = 11
$ws.Range("J401").Value = 11
$ws.Range("K401").Value = 0
$ws.Range("L401").Value = 0
$ws.Range("M401").Value = 0
$ws.Range("N401").Value = 46
$ws.Range("O401").Value = 0
$ws.Range("P401").Value = 0
$ws.Range("Q401").Value = 0

$ws.Range("A402").Value = 45614
$ws.Range("A402").NumberFormat = $dateFmt
$ws.Range("B402").Value = 3948
$ws.Range("C402").Value = 3948
$ws.Range("D402").Value = 3564
$ws.Range("E402").Value = 3613.64990234375
$ws.Range("G402").Value = 3344437
$ws.Range("H402").Value = 2024
$ws.Range("I402").Value = 11
$ws.Range("J402").Value = 18
$ws.Range("K402").Value = 0
$ws.Range("L402").Value = 0
$ws.Range("M402").Value = 0
$ws.Range("N402").Value = 47
$ws.Range("O402").Value = 0
$ws.Range("P402").Value = 0
$ws.Range("Q402").Value = 0

$ws.Range("A403").Value = 45621
$ws.Range("A403").NumberFormat = $dateFmt
$ws.Range("B403").Value = 3686
$ws.Range("C403").Value = 3757.89990234375
$ws.Range("D403").Value = 3592.10009765625
$ws.Range("E403").Value = 3709.60009765625
$ws.Range("G403").Value = 3521889
$ws.Range("H403").Value = 2024
$ws.Range("I403").Value = 11
$ws.Range("J403").Value = 25
$ws.Range("K403").Value = 0
$ws.Range("L403").Value = 0
$ws.Range("M403").Value = 0
$ws.Range("N403").Value = 48
$ws.Range("O403").Value = 0
$ws.Range("P403").Value = 0
$ws.Range("Q403").Value = 0

$ws.Range("A404").Value = 45628
$ws.Range("A404").NumberFormat = $dateFmt
$ws.Range("B404").Value = 3715.449951171875
$ws.Range("C404").Value = 3920.449951171875
$ws.Range("D404").Value = 3622.699951171875
$ws.Range("E404").Value = 3805.550048828125
$ws.Range("G404").Value = 4356276
$ws.Range("H404").Value = 2024
$ws.Range("I404").Value = 12
$ws.Range("J404").Value = 2
$ws.Range("K404").Value = 0
$ws.Range("L404").Value = 0
$ws.Range("M404").Value = 0
$ws.Range("N404").Value = 49
$ws.Range("O404").Value = 0
$ws.Range("P404").Value = 0
$ws.Range("Q404").Value = 0

$ws.Range("A405").Value = 45635
$ws.Range("A405").NumberFormat = $dateFmt
$ws.Range("B405").Value = 3800
$ws.Range("C405").Value = 3853
$ws.Range("D405").Value = 3612
$ws.Range("E405").Value = 3652.300048828125
$ws.Range("G405").Value = 4744458
$ws.Range("H405").Value = 2024
$ws.Range("I405").Value = 12
$ws.Range("J405").Value = 9
$ws.Range("K405").Value = 0
$ws.Range("L405").Value = 0
$ws.Range("M405").Value = 0
$ws.Range("N405").Value = 50
$ws.Range("O405").Value = 0
$ws.Range("P405").Value = 0
$ws.Range("Q405").Value = 0

$ws.Range("A406").Value = 45642
$ws.Range("A406").NumberFormat = $dateFmt
$ws.Range("B406").Value = 3650
$ws.Range("C406").Value = 3672.5
$ws.Range("D406").Value = 3399
$ws.Range("E406").Value = 3408.300048828125
$ws.Range("G406").Value = 3932008
$ws.Range("H406").Value = 2024
$ws.Range("I406").Value = 12
$ws.Range("J406").Value = 16
$ws.Range("K406").Value = 0
$ws.Range("L406").Value = 0
$ws.Range("M406").Value = 0
$ws.Range("N406").Value = 51
$ws.Range("O406").Value = 2
$ws.Range("P406").Value = 0
$ws.Range("Q406").Value = 0

$ws.Range("A407").Value = 45649
$ws.Range("A407").NumberFormat = $dateFmt
$ws.Range("B407").Value = 3443.449951171875
$ws.Range("C407").Value = 3573.39990234375
$ws.Range("D407").Value = 3419.050048828125
$ws.Range("E407").Value = 3568.35009765625
$ws.Range("G407").Value = 2170591
$ws.Range("H407").Value = 2024
$ws.Range("I407").Value = 12
$ws.Range("J407").Value = 23
$ws.Range("K407").Value = 0
$ws.Range("L407").Value = 0
$ws.Range("M407").Value = 0
$ws.Range("N407").Value = 52
$ws.Range("O407").Value = 0
$ws.Range("P407").Value = 0
$ws.Range("Q407").Value = 0

$ws.Range("A408").Value = 45656
$ws.Range("A408").NumberFormat = $dateFmt
$ws.Range("B408").Value = 3570
$ws.Range("C408").Value = 4165.89990234375
$ws.Range("D408").Value = 3495.050048828125
$ws.Range("E408").Value = 4025.199951171875
$ws.Range("G408").Value = 13133195
$ws.Range("H408").Value = 2024
$ws.Range("I408").Value = 12
$ws.Range("J408").Value = 30
$ws.Range("K408").Value = 0
$ws.Range("L408").Value = 0
$ws.Range("M408").Value = 0
$ws.Range("N408").Value = 1
$ws.Range("O408").Value = 0
$ws.Range("P408").Value = 0
$ws.Range("Q408").Value = 0

$ws.Range("A409").Value = 45663
$ws.Range("A409").NumberFormat = $dateFmt
$ws.Range("B409").Value = 4025.199951171875
$ws.Range("C409").Value = 4061.75
$ws.Range("D409").Value = 3666.64990234375
$ws.Range("E409").Value = 3686.25
$ws.Range("G409").Value = 4040905
$ws.Range("H409").Value = 2025
$ws.Range("I409").Value = 1
$ws.Range("J409").Value = 6
$ws.Range("K409").Value = 0
$ws.Range("L409").Value = 0
$ws.Range("M409").Value = 0
$ws.Range("N409").Value = 2
$ws.Range("O409").Value = 0
$ws.Range("P409").Value = 0
$ws.Range("Q409").Value = 0

$ws.Range("A410").Value = 45670
$ws.Range("A410").NumberFormat = $dateFmt
$ws.Range("B410").Value = 3501
$ws.Range("C410").Value = 3649.949951171875
$ws.Range("D410").Value = 3443
$ws.Range("E410").Value = 3620.64990234375
$ws.Range("G410").Value = 5083577
$ws.Range("H410").Value = 2025
$ws.Range("I410").Value = 1
$ws.Range("J410").Value = 13
$ws.Range("K410").Value = 0
$ws.Range("L410").Value = 0
$ws.Range("M410").Value = 0
$ws.Range("N410").Value = 3
$ws.Range("O410").Value = 0
$ws.Range("P410").Value = 0
$ws.Range("Q410").Value = 0

$ws.Range("A411").Value = 45677
$ws.Range("A411").NumberFormat = $dateFmt
$ws.Range("B411").Value = 3624
$ws.Range("C411").Value = 3640
$ws.Range("D411").Value = 3524.5
$ws.Range("E411").Value = 3579.949951171875
$ws.Range("G411").Value = 1280121
$ws.Range("H411").Value = 2025
$ws.Range("I411").Value = 1
$ws.Range("J411").Value = 20
$ws.Range("K411").Value = 0
$ws.Range("L411").Value = 0
$ws.Range("M411").Value = 0
$ws.Range("N411").Value = 4
$ws.Range("O411").Value = 0
$ws.Range("P411").Value = 0
$ws.Range("Q411").Value = 0

